$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1405
$ws1.Range("F5").Value = 320
$ws1.Range("F6").Value = 1036
$ws1.Range("F7").Value = 10757
$ws1.Range("F11").Value = 1040
$ws1.Range("F13").Value = 12077
$ws1.Range("F14").Value = 12537

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1405
$ws4.Range("F6").Value = 320
$ws4.Range("F7").Value = 1036
$ws4.Range("F8").Value = 10757
$ws4.Range("F12").Value = 1040
$ws4.Range("F14").Value = 12077
$ws4.Range("F15").Value = 12537
